$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Отобранные участники")

# ---------------------------------------------------------------------------
# 1. Copy formatting for the 4 new rows (244-247) from existing rows that
#    already carry the exact style combination we need:
#      - row 211 -> rows where column B is a hyperlink (style 20)
#      - row 225 -> rows where column B is plain text (style 13)
# ---------------------------------------------------------------------------
$ws.Range("A211:U211").Copy()
$ws.Range("A244:U244").PasteSpecial(-4122)

$ws.Range("A225:U225").Copy()
$ws.Range("A245:U245").PasteSpecial(-4122)

$ws.Range("A211:U211").Copy()
$ws.Range("A246:U246").PasteSpecial(-4122)

$ws.Range("A225:U225").Copy()
$ws.Range("A247:U247").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Row 244 - Чебыкина Юлия Александровна
# ---------------------------------------------------------------------------
$ws.Cells.Item(244,1).Value = "Чебыкина Юлия Александровна"
$ws.Cells.Item(244,2).Value = "https://t.me/YaMusenka"
$ws.Cells.Item(244,3).Value = 89091011335
$ws.Cells.Item(244,4).Value = "Android"
$ws.Cells.Item(244,5).Value = "Пермский край"
$ws.Cells.Item(244,6).Value = "Пермь"
$ws.Cells.Item(244,7).Value = "Женский"
$ws.Cells.Item(244,8).Value = 40
$ws.Cells.Item(244,10).Value = "Нет, не являюсь клиентом данного банка"
$ws.Cells.Item(244,11).Value = "Да"
$ws.Cells.Item(244,12).Value = "Запас"
$ws.Cells.Item(244,13).Value = "Нет, не являюсь клиентом данного банка"
$ws.Cells.Item(244,16).Value = "Без обращения"
$ws.Cells.Item(244,17).Value = "Сценарий 4"

# ---------------------------------------------------------------------------
# 3. Row 245 - Булатова Елена Анатольевна
# ---------------------------------------------------------------------------
$ws.Cells.Item(245,1).Value = "Булатова Елена Анатольевна"
$ws.Cells.Item(245,2).Value = "@elenabulatti"
$ws.Cells.Item(245,3).Value = 89276561018
$ws.Cells.Item(245,4).Value = "Android"
$ws.Cells.Item(245,5).Value = "Самарская область "
$ws.Cells.Item(245,6).Value = "Самара"
$ws.Cells.Item(245,7).Value = "Женский"
$ws.Cells.Item(245,8).Value = 39
$ws.Cells.Item(245,10).Value = "Нет, не являюсь клиентом данного банка"
$ws.Cells.Item(245,11).Value = "Да"
$ws.Cells.Item(245,12).Value = "Нет, не являюсь клиентом данного банка"
$ws.Cells.Item(245,13).Value = "Да"
$ws.Cells.Item(245,16).Value = "Сценарий 1"
$ws.Cells.Item(245,17).Value = "Сценарий 4"
$ws.Cells.Item(245,20).Value = "Сценарий 2"
$ws.Cells.Item(245,21).Value = "Без обращения"

# ---------------------------------------------------------------------------
# 4. Row 246 - Терещенко Юрий Андреевич
# ---------------------------------------------------------------------------
$ws.Cells.Item(246,1).Value = "Терещенко Юрий Андреевич"
$ws.Cells.Item(246,2).Value = "t.me/uriyter "
$ws.Cells.Item(246,3).Value = 89896364327
$ws.Cells.Item(246,4).Value = "Android"
$ws.Cells.Item(246,5).Value = "Ростов-на-Дону "
$ws.Cells.Item(246,6).Value = "Ростов-на-Дону "
$ws.Cells.Item(246,7).Value = "Мужской"
$ws.Cells.Item(246,8).Value = 29
$ws.Cells.Item(246,10).Value = "Запас"
$ws.Cells.Item(246,11).Value = "Да"
$ws.Cells.Item(246,12).Value = "Запас"
$ws.Cells.Item(246,13).Value = "Нет, не являюсь клиентом данного банка"
$ws.Cells.Item(246,16).Value = "Без обращения"
$ws.Cells.Item(246,17).Value = "Сценарий 4"

# ---------------------------------------------------------------------------
# 5. Row 247 - Дмитриева Олеся Владимировна
# ---------------------------------------------------------------------------
$ws.Cells.Item(247,1).Value = "Дмитриева Олеся Владимировна"
$ws.Cells.Item(247,2).Value = "Olani11"
$ws.Cells.Item(247,3).Value = 89045883559
$ws.Cells.Item(247,4).Value = "Android"
$ws.Cells.Item(247,5).Value = "Омск"
$ws.Cells.Item(247,6).Value = "Омск"
$ws.Cells.Item(247,7).Value = "Женский"
$ws.Cells.Item(247,8).Value = 39
$ws.Cells.Item(247,10).Value = "Запас"
$ws.Cells.Item(247,11).Value = "Да"
$ws.Cells.Item(247,12).Value = "Да"
$ws.Cells.Item(247,13).Value = "Да"
$ws.Cells.Item(247,16).Value = "Сценарий 1"
$ws.Cells.Item(247,17).Value = "Сценарий 4"
$ws.Cells.Item(247,18).Value = "Сценарий 2"
$ws.Cells.Item(247,19).Value = "Сценарий 5"
$ws.Cells.Item(247,20).Value = "Сценарий 2"
$ws.Cells.Item(247,21).Value = "Сценарий 5"

# ---------------------------------------------------------------------------
# 6. Hyperlinks on column B for the two Telegram-link rows (244 and 246),
#    then restore the sheet's usual hyperlink look (style copied from B211)
#    since Hyperlinks.Add applies Excel's built-in hyperlink style.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Cells.Item(244,2), "https://t.me/YaMusenka")
$ws.Hyperlinks.Add($ws.Cells.Item(246,2), "t.me/uriyter ")

$ws.Range("B211").Copy()
$ws.Range("B244").PasteSpecial(-4122)
$ws.Range("B211").Copy()
$ws.Range("B246").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 7. Cell-level corrections on pre-existing rows
# ---------------------------------------------------------------------------
$ws.Cells.Item(238,18).Value = "Сценарий 1"
$ws.Cells.Item(240,10).Value = "Да"
$ws.Cells.Item(240,14).Value = "Сценарий 2"
$ws.Cells.Item(240,15).Value = "Без обращения"
$ws.Cells.Item(241,20).Value = "Сценарий 1"
$ws.Cells.Item(241,21).Value = "Без обращения"

# ---------------------------------------------------------------------------
# 8. Extend the AutoFilter range to cover the new rows
# ---------------------------------------------------------------------------
$ws.Range("A1:AA243").AutoFilter()
$ws.Range("A1:AA247").AutoFilter()

# ---------------------------------------------------------------------------
# 9. Update the _xlnm._FilterDatabase defined name to match the new range
# ---------------------------------------------------------------------------
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='Отобранные участники'!`$A`$1:`$AA`$247"

# ---------------------------------------------------------------------------
# 10. Recalculate so cached formula results (COUNTIF helpers) are refreshed
# ---------------------------------------------------------------------------
$excel.CalculateFullRebuild()
